# ABC-analysis columns (Аккум.доля / Категория) added to the existing
# "Доля" table, plus swap of the Товар 3 / Товар 4 sales figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix rows 4/5: sales (C) and Доля (D) values were swapped -------------
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 10
$ws.Range("C5").Value = 8
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "4.371585"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null

# --- new headers (E1/F1) ---------------------------------------------------
$ws.Range("E1").Value = "Аккум.доля"
$ws.Range("F1").Value = "Категория"

# --- helper cell used later to reset formatting to the sheet's default style
# (D1 already carries the "no explicit style" look we want to replicate)

# --- column E: cumulative share (Аккум.доля) -------------------------------
# E2 repeats the same text value as D2 ("54.644809")
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "54.644809"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "81.967213"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "92.896175"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "97.267760"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null

# E6 is a genuine number (100), not text
$ws.Range("E6").Value = 100

# --- column F: ABC category, formatted like the rest of the data rows -----
$ws.Range("C2:C6").Copy() | Out-Null
$ws.Range("F2:F6").PasteSpecial(-4122) | Out-Null

$ws.Range("F2").Value = "A"
$ws.Range("F3").Value = "B"
$ws.Range("F4").Value = "B"
$ws.Range("F5").Value = "C"
$ws.Range("F6").Value = "C"

# --- column width for the new "Аккум.доля" column --------------------------
$ws.Columns.Item(5).ColumnWidth = 11.25

# --- selection / active cell as left by the author -------------------------
$ws.Range("F12").Select() | Out-Null
